$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 14
$ws.Cells.Item($row, 1).Value = 42619.891273148147
$ws.Cells.Item($row, 2).Value = -8
$ws.Cells.Item($row, 3).Value = 50
$ws.Cells.Item($row, 4).Value = 47
$ws.Cells.Item($row, 5).Value = 50
$ws.Cells.Item($row, 6).Value = 71
$ws.Cells.Item($row, 7).Value = 10007
$ws.Cells.Item($row, 8).Value = 15607
$ws.Cells.Item($row, 9).Value = 1691
$ws.Cells.Item($row, 10).Value = 187
$ws.Cells.Item($row, 11).Value = 175
$ws.Cells.Item($row, 12).Value = 4
$ws.Cells.Item($row, 13).Value = 10
$ws.Cells.Item($row, 14).Value = "Bag"
